$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Nlgn3"
$ws.Cells.Item(2, 3).Value = "Nrxn1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3239616666666667
$ws.Cells.Item(2, 8).Value = 0.9718850000000001
$ws.Cells.Item(2, 9).Value = 0.1555281019885789
$ws.Cells.Item(2, 10).Value = 0.1555281019885789
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.07823633333333334
$ws.Cells.Item(2, 14).Value = 0.234709
$ws.Cells.Item(2, 15).Value = 0.07325462794193288
$ws.Cells.Item(2, 16).Value = 0.07325462794193287
$ws.Cells.Item(2, 17).Value = 0.02534557294055556
$ws.Cells.Item(2, 18).Value = 0.228110156465
$ws.Cells.Item(2, 19).Value = 0.01139315324568834
$ws.Cells.Item(2, 20).Value = 0.01139315324568833

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Nlgn3"
$ws.Cells.Item(3, 3).Value = "Nrxn1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3239616666666667
$ws.Cells.Item(3, 8).Value = 0.9718850000000001
$ws.Cells.Item(3, 9).Value = 0.1555281019885789
$ws.Cells.Item(3, 10).Value = 0.1555281019885789
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.989769
$ws.Cells.Item(3, 14).Value = 2.969307
$ws.Cells.Item(3, 15).Value = 0.9267453720580672
$ws.Cells.Item(3, 16).Value = 0.9267453720580671
$ws.Cells.Item(3, 17).Value = 0.320647214855
$ws.Cells.Item(3, 18).Value = 2.885824933695
$ws.Cells.Item(3, 19).Value = 0.1441349487428906
$ws.Cells.Item(3, 20).Value = 0.1441349487428905

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Nlgn3"
$ws.Cells.Item(4, 3).Value = "Nrxn1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.026996666666667
$ws.Cells.Item(4, 8).Value = 3.08099
$ws.Cells.Item(4, 9).Value = 0.4930424144274184
$ws.Cells.Item(4, 10).Value = 0.4930424144274184
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.07823633333333334
$ws.Cells.Item(4, 14).Value = 0.234709
$ws.Cells.Item(4, 15).Value = 0.07325462794193288
$ws.Cells.Item(4, 16).Value = 0.07325462794193287
$ws.Cells.Item(4, 17).Value = 0.08034845354555555
$ws.Cells.Item(4, 18).Value = 0.72313608191
$ws.Cells.Item(4, 19).Value = 0.03611763862847282
$ws.Cells.Item(4, 20).Value = 0.03611763862847281

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Nlgn3"
$ws.Cells.Item(5, 3).Value = "Nrxn1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.026996666666667
$ws.Cells.Item(5, 8).Value = 3.08099
$ws.Cells.Item(5, 9).Value = 0.4930424144274184
$ws.Cells.Item(5, 10).Value = 0.4930424144274184
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.989769
$ws.Cells.Item(5, 14).Value = 2.969307
$ws.Cells.Item(5, 15).Value = 0.9267453720580672
$ws.Cells.Item(5, 16).Value = 0.9267453720580671
$ws.Cells.Item(5, 17).Value = 1.01648946377
$ws.Cells.Item(5, 18).Value = 9.14840517393
$ws.Cells.Item(5, 19).Value = 0.4569247757989456
$ws.Cells.Item(5, 20).Value = 0.4569247757989455

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Nlgn3"
$ws.Cells.Item(6, 3).Value = "Nrxn1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.73202
$ws.Cells.Item(6, 8).Value = 2.19606
$ws.Cells.Item(6, 9).Value = 0.3514294835840027
$ws.Cells.Item(6, 10).Value = 0.3514294835840027
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.07823633333333334
$ws.Cells.Item(6, 14).Value = 0.234709
$ws.Cells.Item(6, 15).Value = 0.07325462794193288
$ws.Cells.Item(6, 16).Value = 0.07325462794193287
$ws.Cells.Item(6, 17).Value = 0.05727056072666667
$ws.Cells.Item(6, 18).Value = 0.5154350465400001
$ws.Cells.Item(6, 19).Value = 0.02574383606777173
$ws.Cells.Item(6, 20).Value = 0.02574383606777172

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Nlgn3"
$ws.Cells.Item(7, 3).Value = "Nrxn1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.73202
$ws.Cells.Item(7, 8).Value = 2.19606
$ws.Cells.Item(7, 9).Value = 0.3514294835840027
$ws.Cells.Item(7, 10).Value = 0.3514294835840027
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.989769
$ws.Cells.Item(7, 14).Value = 2.969307
$ws.Cells.Item(7, 15).Value = 0.9267453720580672
$ws.Cells.Item(7, 16).Value = 0.9267453720580671
$ws.Cells.Item(7, 17).Value = 0.72453070338
$ws.Cells.Item(7, 18).Value = 6.52077633042
$ws.Cells.Item(7, 19).Value = 0.325685647516231
$ws.Cells.Item(7, 20).Value = 0.325685647516231
